$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For data rows 2-4, the "phase" column (J) value is removed and subsequent
# columns (K, L, ...) shift right by one, with 0 inserted at J.
for ($r = 2; $r -le 4; $r++) {
    $oldJ = $ws.Cells.Item($r, 10).Value2  # column J
    $oldK = $ws.Cells.Item($r, 11).Value2  # column K
    $oldL = $ws.Cells.Item($r, 12).Value2  # column L

    $ws.Cells.Item($r, 10).Value2 = 0
    $ws.Cells.Item($r, 11).Value2 = $oldJ
    $ws.Cells.Item($r, 12).Value2 = $oldK
    $ws.Cells.Item($r, 13).Value2 = $oldL
}

# Update latitude (column X) precision for rows 2-4.
$ws.Range("X2").Value2 = 34.1321878099
$ws.Range("X3").Value2 = 34.1321878099
$ws.Range("X4").Value2 = 34.1321878099
